$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10759000
$ws.Range("C3").Value = -539999.9999999995
$ws.Range("C4").Value = 3096428.571428571
$ws.Range("C5").Value = 2060000
$ws.Range("C6").Value = 6980000
$ws.Range("C7").Value = 810000
$ws.Range("C8").Value = 6000000
$ws.Range("C9").Value = 16597500
$ws.Range("C10").Value = 17189333.33333333
$ws.Range("C11").Value = 19611571.42857143
$ws.Range("C12").Value = 15271000
$ws.Range("C13").Value = 97834833.33333333
